$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.005.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5138"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3835"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08266"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.110"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.198"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.27%  "

$ws.Range("E13").Value = "  -0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.866.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.255"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06655"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.86%  "

$ws.Range("E20").Value = "  -2.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("E22").Value = "  -1.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.036.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.253"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.075.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.518"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.70%  "

$ws.Range("E31").Value = "  +1.12%  "

$ws.Range("E32").Value = "  -3.08%  "

$ws.Range("E33").Value = "  +6.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.604"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.357"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.66%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06496"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6567"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.64%  "

$ws.Range("E40").Value = "  -0.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.008"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.226"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.88%  "

$ws.Range("E44").Value = "  +2.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.659"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("E48").Value = "  +0.66%  "

$ws.Range("E49").Value = "  -2.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("E51").Value = "  -1.55%  "
